$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.648.60'
$ws.Range('E2').Value = '  +2.96%  '
$ws.Range('D3').Value = '2.190.10'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '81.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.70%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0920'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.59%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = '2.518.41'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.31'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '2.176.96'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.778'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').Value = '43.597.56'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000103'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +12.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.74%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '41.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.22'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0876'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.113'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0352'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +15.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.10'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '62.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.73%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0985'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('E48').Value = '  +3.92%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +26.91%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.437'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.55%  '
